# sentence calculator version 1
# calculate the length of the sentences in each book
#
# Insert a new header row at the top of the chapter list and label it
# "chapters". This shifts every existing chapter-name row down by one.
# Also remove the stray "诺伯" note that lived in column D next to
# "第14章" (now shifted to column D of the row holding "第14章").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push all existing rows down by one and create the new header row.
$ws.Rows("1:1").Insert()
$ws.Range("A1").Value = "chapters"

# The old "诺伯" annotation (originally D14) shifted down to D15 along
# with the rest of the data; drop it, leaving only the single chapters
# column behind.
$ws.Range("D15").ClearContents()

# Restore the default selection at the top-left cell.
$ws.Range("A1").Select()
